# Update predicted labels for model v1 - improved results
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "B5"  = 1
    "B9"  = 1
    "B10" = 0
    "B11" = 1
    "B12" = 1
    "B17" = 0
    "B18" = 0
    "B22" = 1
    "B24" = 1
    "B28" = 0
    "B40" = 0
    "B48" = 1
    "B52" = 0
    "B55" = 0
    "B57" = 1
    "B60" = 1
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
